$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.891.20"
$ws.Range("E2").Value = "  +0.81%  "
$ws.Range("D3").Value = "2.521.21"
$ws.Range("E3").Value = "  -1.99%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.14%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.23%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.79%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.526"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.52%  "
$ws.Range("D9").Value = "2.519.46"
$ws.Range("E9").Value = "  -2.09%  "
$ws.Range("E10").Value = "  -0.07%  "
$ws.Range("E11").Value = "  +1.90%  "
$ws.Range("E12").Value = "  -0.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.342"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.55%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "26.59"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "3.005.10"
$ws.Range("E15").Value = "  -1.35%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000175"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.10%  "
$ws.Range("D17").Value = "66.789.42"
$ws.Range("E17").Value = "  +0.11%  "
$ws.Range("D18").Value = "2.490.82"
$ws.Range("E18").Value = "  -3.87%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.97"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.30"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.09%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "353.37"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.97%  "
$ws.Range("E22").Value = "  -1.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.33%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.97"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.44%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "69.55"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.40%  "
$ws.Range("D28").Value = "2.679.21"
$ws.Range("E28").Value = "  -1.08%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("D30").Value = "0.0₃0977"
$ws.Range("E30").Value = "  -0.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "530.73"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.08%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.45%  "
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.67%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.999"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("B37").Value = "Monero"
$ws.Range("C37").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "157.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.58%  "
$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.54%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.53"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("E41").Value = "  -1.54%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.08"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.33%  "
$ws.Range("E44").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +4.93%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "39.62"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.01%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "148.54"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.34%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.555"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("E49").Value = "  -3.39%  "
$ws.Range("E50").Value = "  -1.25%  "
$ws.Range("E51").Value = "  +0.40%  "
